$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, bordered, centered) from an existing header cell (A1)
# onto the new header cells, then set their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Cells.Item(1, 30).Value2 = "Wins"
$ws.Cells.Item(1, 31).Value2 = "Losses"
$ws.Cells.Item(1, 32).Value2 = "Ties"

# Fill in the team record (Wins/Losses/Ties) for every player data row.
for ($r = 2; $r -le 55; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 57
    $ws.Cells.Item($r, 31).Value2 = 105
    $ws.Cells.Item($r, 32).Value2 = 0
}
